$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The second-ranked player (row 3, "Remco Riem") was removed from the
# standings. Deleting the whole row shifts every row below it up by one,
# which matches rows 4-22 of the old sheet becoming rows 3-21 of the new
# sheet, and updates the sheet dimension from A1:J22 to A1:J21 automatically.
$ws.Rows("3").Delete()

# Column A ("Rang") holds the competitor's rank (ties share a rank), derived
# from column I ("Totaal"). With one player removed the ranks below the gap
# compress, so recompute the cached rank values for the remaining 20 players
# (now in rows 2-21).
$rangValues = @(1, 2, 3, 3, 3, 6, 7, 7, 9, 10, 11, 12, 12, 12, 15, 16, 17, 17, 19, 20)

for ($i = 0; $i -lt $rangValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $rangValues[$i]
}
